# Apply the "parte 1 de nuevos estado de cuenta" update to the EC workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared "Periodo Mora" label used by every worker row (2507 -> 2508)
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2508"
$ws.Range("E19").Value = "2508"
$ws.Range("E20").Value = "2508"

# Total "Valor Mora" at the top of the statement
$ws.Range("E11").Value = 407760

# Per-worker "Valor Mora" (F) and "Salario Basico" (G) figures
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

$ws.Range("F18").Value = 180000
$ws.Range("G18").Value = 4500000

$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500
